$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Estadisticos 2P" - update the stat columns for all 4 groups
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

# Row 2 (6ALCM)
$ws2.Range("D2").Value = 7
$ws2.Range("E2").Value = 10
$ws2.Range("F2").Value = 26
$ws2.Range("G2").Value = 72.22
$ws2.Range("H2").Value = 7.6

# Row 3 (6APM)
$ws2.Range("D3").Value = 7
$ws2.Range("E3").Value = 14
$ws2.Range("F3").Value = 8
$ws2.Range("G3").Value = 36.36
$ws2.Range("H3").Value = 6.3

# Row 4 (6ARHM)
$ws2.Range("D4").Value = 9
$ws2.Range("E4").Value = 9
$ws2.Range("F4").Value = 25
$ws2.Range("G4").Value = 73.53
$ws2.Range("H4").Value = 8.9

# Row 5 (6BLCM)
$ws2.Range("D5").Value = 18
$ws2.Range("E5").Value = 20
$ws2.Range("F5").Value = 15
$ws2.Range("G5").Value = 42.86
$ws2.Range("H5").Value = 7.1

# -----------------------------------------------------------------
# Sheet "Estadisticos Final" - row 5 (6BLCM) stats tweak
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("E5").Value = 17
$ws3.Range("F5").Value = 18
$ws3.Range("G5").Value = 51.43
$ws3.Range("H5").Value = 6.4

# -----------------------------------------------------------------
# Sheet "Rescatables" - add the 14 rescued students
# Data is entered column-by-column (Paterno, then Materno, then
# Nombres) to mirror how the surnames/names were typed in originally.
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$claves = @(18330051920261,18330051920302,18330051920308,18330051920311,18330051920440,18330051920366,18330051920217,18330051920329,18330051920424,18330051920333,18330051920339,18330051920342,18330051920347,18330051920352)
$paterno = @("DE GANTE","SANCHEZ","VALDES","VASQUEZ","DURAN","GARCIA","MARQUEZ","GUERRA","LOPEZ","LOPEZ","MORALES","REYES","TEXOCO","ZEPAHUA")
$materno = @("GUTIERREZ","REYES","MARIN","PIEDRAS","CORTES","ARENAS","HERNANDEZ","ROMERO","APALE","MORALES","TREJO","SARMIENTO","DE JESUS","JUAREZ")
$nombres = @("DANIELA","YAIRA GUADALUPE","BRANDON","MARTHA MARISOL","LUIS GABRIEL","JOSUA","ANDRES","JOCELYN","MARIA LARET","SILVIA ESMERALDA","ROCIO TAMARA","INGRID PAOLA","MAYTE","QUETZALI")
$nombreLargo = @("TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA","TEMAS DE FILOSOFÍA")
$grupo = @("6ALCM","6ALCM","6ALCM","6ALCM","6APM","6APM","6ARHM","6BLCM","6BLCM","6BLCM","6BLCM","6BLCM","6BLCM","6BLCM")
$reprobadas = @(2,2,2,2,2,2,2,2,2,2,2,2,2,2)

for ($i = 0; $i -lt $paterno.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 2).Value = $paterno[$i]
}
for ($i = 0; $i -lt $materno.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 3).Value = $materno[$i]
}
for ($i = 0; $i -lt $nombres.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 4).Value = $nombres[$i]
}
for ($i = 0; $i -lt $claves.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 1).Value = $claves[$i]
    $ws4.Cells.Item($row, 5).Value = $nombreLargo[$i]
    $ws4.Cells.Item($row, 6).Value = $grupo[$i]
    $ws4.Cells.Item($row, 7).Value = $reprobadas[$i]
}
